# Apply averaged-intensity results including the new spiral sampling schemes
# (Gaussian-Quadrature row relocated, 3 new Spiral-* rows inserted, remaining rows shifted down)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=10; A=8; B="Gaussian-Quadrature"; Vals=@(1.365148640725282, 0.3837816257825609, 1.051021841029485, 1.365148640725282, 0.6851569995871816, 1.134651966991828, 1.133176527350012, 0.3837816257825609, 0.7174017334060228, 1.041275187065652, 0.958822933577725) },
    @{ Row=11; A=9; B="Spiral-90deg-10rot-5space"; Vals=@(0.7745066664228107, 1.240672328105504, 1.112867198188798, 0.7745066664228107, 0.8451038369145706, 1.633819030898185, 0.9916895547144543, 1.240672328105504, 1.176769763147151, 0.9756382147849808, 1.099776435874054) },
    @{ Row=12; A=10; B="Spiral-90deg-15rot-5space"; Vals=@(0.7758502712701368, 1.242654337317725, 1.111644851960693, 0.7758502712701368, 0.8458514540075714, 1.631980417083467, 0.9909160839665972, 1.242654337317725, 1.177149594639209, 0.9764999329546727, 1.099816235934365) },
    @{ Row=13; A=11; B="Spiral-90deg-10rot-3space"; Vals=@(0.7752097379639964, 1.240158767427029, 1.113078474216721, 0.7752097379639964, 0.844993179444269, 1.631513064206297, 0.991650365349732, 1.240158767427029, 1.176618620821875, 0.9759141793929357, 1.099433931434674) },
    @{ Row=14; A=12; B="NoRotation-tilt60deg"; Vals=@(0.5927680000000006, 0.5840240000000011, 1.309084, 0.5927680000000006, 0.3867640000000004, 3.362899999999992, 1.143523999999999, 0.5840240000000011, 0.9465540000000008, 0.7696610000000006, 1.229843999999999) },
    @{ Row=15; A=13; B="Rotation-NoTilt"; Vals=@(0.21, 0, 1.551112499999998, 0.21, 0.01, 5.087262499999993, 1.2698625, 0, 0.7755562499999992, 0.4927781249999996, 1.354706249999999) },
    @{ Row=16; A=14; B="Rotation-60detTilt"; Vals=@(0.5507476966400019, 0.4137395924992008, 1.327566295449597, 0.5507476966400019, 0.4147668774912002, 3.308730519449605, 1.155351315660798, 0.4137395924992008, 0.8706529439743989, 0.7107003203072004, 1.195150382865067) },
    @{ Row=17; A=15; B="HexGrid-90degTilt5degRes"; Vals=@(0.9792400932042136, 0.9942469434978274, 0.9981685128636005, 0.9792400932042136, 0.9884254592161964, 1.001347330476415, 0.9935757308028919, 0.9942469434978274, 0.996207728180714, 0.9877239106924638, 0.9925006783435242) },
    @{ Row=18; A=16; B="HexGrid-90degTilt22p5degRes"; Vals=@(1.067063351873834, 1.07691332059648, 0.9326994561095385, 1.067063351873834, 1.02620181448143, 0.9658785909203682, 0.9550665704039035, 1.07691332059648, 1.004806388353009, 1.035934870113422, 1.003970517397593) },
    @{ Row=19; A=17; B="HexGrid-60degTilt5degRes"; Vals=@(0.9768844703237514, 1.233286640664121, 0.9286378950768857, 0.9768844703237514, 1.114010473699339, 0.777109948640353, 0.9326324050533024, 1.233286640664121, 1.080962267870503, 1.028923369097127, 0.9937603055762919) }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $col = 3
    foreach ($v in $r.Vals) {
        $ws.Cells.Item($r.Row, $col).Value = $v
        $col = $col + 1
    }
}

# Rows 17-19 are brand new rows beyond the original sheet extent (A1:M16); give column A
# the same "HKL index" formatting (bold, centered, thin border) used throughout column A by
# copying the format already applied to A16.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
